# Fixed Email field name in examples
# The "Email Message" header in the example data was actually the column
# for an email address, so rename the header to simply "Email".

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range("I1").Value = "Email"
